$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Cells.Item(2, 1).Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f45a99f0580>),`n                ('model',`n                 RandomForestClassifier(max_depth=5, min_samples_leaf=6,`n                                        min_samples_split=5, n_estimators=50,`n                                        random_state=42))])"
$ws.Cells.Item(2, 2).Value = 0.6857142857142857
$ws.Cells.Item(2, 3).Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a64541040>, 'scaler': None, 'model__n_estimators': 50, 'model__min_samples_split': 5, 'model__min_samples_leaf': 6, 'model__max_features': 'sqrt', 'model__max_depth': 5}"
$ws.Cells.Item(2, 4).Value = 0.3333333333333333
$ws.Cells.Item(2, 5).Value = "[1 1 0 0 1 0 0 0 0 1 0 1]"
$ws.Cells.Item(2, 6).Value = "[0 1 1 0 0 1 1 1 1 1 0 0]"
$ws.Cells.Item(2, 7).Value = 77
$ws.Cells.Item(2, 8).Value = 0.8175406871609404
$ws.Cells.Item(2, 9).Value = 0.02377518512520467
$ws.Cells.Item(2, 10).Value = 0.5918022905364677
$ws.Cells.Item(2, 11).Value = 0.07405076486267803

# ---- Row 3 ----
$ws.Cells.Item(3, 1).Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f3a6464e8b0>),`n                ('model',`n                 RandomForestClassifier(max_depth=2, max_features='log2',`n                                        min_samples_leaf=5, min_samples_split=4,`n                                        n_estimators=50, random_state=42))])"
$ws.Cells.Item(3, 2).Value = 0.6666666666666666
$ws.Cells.Item(3, 3).Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a64619520>, 'scaler': None, 'model__n_estimators': 50, 'model__min_samples_split': 4, 'model__min_samples_leaf': 5, 'model__max_features': 'log2', 'model__max_depth': 2}"
$ws.Cells.Item(3, 4).Value = 0.7777777777777778
$ws.Cells.Item(3, 5).Value = "[1 1 0 1 0 0 1 0 1 1 1 0]"
$ws.Cells.Item(3, 6).Value = "[1 1 1 1 1 1 1 1 1 1 1 0]"
$ws.Cells.Item(3, 7).Value = 69
$ws.Cells.Item(3, 8).Value = 0.8103571428571428
$ws.Cells.Item(3, 9).Value = 0.02569604584235274
$ws.Cells.Item(3, 10).Value = 0.5394047619047618
$ws.Cells.Item(3, 11).Value = 0.0824882319080132

# ---- Row 4 ----
$ws.Cells.Item(4, 1).Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f3a641576d0>),`n                ('model',`n                 RandomForestClassifier(max_depth=4, max_features='log2',`n                                        min_samples_leaf=5, n_estimators=50,`n                                        random_state=42))])"
$ws.Cells.Item(4, 2).Value = 0.6285714285714284
$ws.Cells.Item(4, 3).Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a646ae340>, 'scaler': None, 'model__n_estimators': 50, 'model__min_samples_split': 2, 'model__min_samples_leaf': 5, 'model__max_features': 'log2', 'model__max_depth': 4}"
$ws.Cells.Item(4, 4).Value = 0.7777777777777777
$ws.Cells.Item(4, 5).Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Cells.Item(4, 6).Value = "[1 1 1 1 1 1 0 0 1 1 1 1]"
$ws.Cells.Item(4, 7).Value = 42
$ws.Cells.Item(4, 8).Value = 0.8326441102756893
$ws.Cells.Item(4, 9).Value = 0.03138730663534699
$ws.Cells.Item(4, 10).Value = 0.5266917293233082
$ws.Cells.Item(4, 11).Value = 0.08374548353770728

# ---- Row 5 (new) ----
$ws.Cells.Item(5, 1).Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f3a64157400>),`n                ('model',`n                 RandomForestClassifier(max_depth=4, max_features='log2',`n                                        min_samples_leaf=5, min_samples_split=3,`n                                        n_estimators=5, random_state=42))])"
$ws.Cells.Item(5, 2).Value = 0.6380952380952382
$ws.Cells.Item(5, 3).Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6464e940>, 'scaler': None, 'model__n_estimators': 5, 'model__min_samples_split': 3, 'model__min_samples_leaf': 5, 'model__max_features': 'log2', 'model__max_depth': 4}"
$ws.Cells.Item(5, 4).Value = 0.6153846153846153
$ws.Cells.Item(5, 5).Value = "[1 1 0 0 0 0 1 0 1 1 1 1]"
$ws.Cells.Item(5, 6).Value = "[0 0 1 0 0 0 1 1 1 0 1 1]"
$ws.Cells.Item(5, 7).Value = 11
$ws.Cells.Item(5, 8).Value = 0.8420024420024419
$ws.Cells.Item(5, 9).Value = 0.02426141396175129
$ws.Cells.Item(5, 10).Value = 0.5230769230769231
$ws.Cells.Item(5, 11).Value = 0.07344992139759915

# ---- Row 6 (new) ----
$ws.Cells.Item(6, 1).Value = "Pipeline(steps=[('scaler', None), ('selector', None),`n                ('model',`n                 RandomForestClassifier(max_depth=2, max_features='log2',`n                                        min_samples_leaf=3, min_samples_split=4,`n                                        n_estimators=10, random_state=42))])"
$ws.Cells.Item(6, 2).Value = 0.6571428571428571
$ws.Cells.Item(6, 3).Value = "{'selector': None, 'scaler': None, 'model__n_estimators': 10, 'model__min_samples_split': 4, 'model__min_samples_leaf': 3, 'model__max_features': 'log2', 'model__max_depth': 2}"
$ws.Cells.Item(6, 4).Value = 0.5714285714285715
$ws.Cells.Item(6, 5).Value = "[1 1 1 1 0 0 0 0 1 1 0 0]"
$ws.Cells.Item(6, 6).Value = "[1 1 1 0 0 1 0 1 1 0 1 1]"
$ws.Cells.Item(6, 7).Value = 14
$ws.Cells.Item(6, 8).Value = 0.845539481615431
$ws.Cells.Item(6, 9).Value = 0.02466981175987563
$ws.Cells.Item(6, 10).Value = 0.5585292344786015
$ws.Cells.Item(6, 11).Value = 0.07476658767319569

Write-Output "edit applied"
